# Applies crypto price/volume updates for Wed Sep 13 2023 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.016.81"
$ws.Cells.Item(3, 4).Value = "1.595.42"
$ws.Cells.Item(3, 5).Value = "  +0.65%  "
$ws.Cells.Item(4, 5).Value = "  -0.04%  "
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "211.53"
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.55%  "
$ws.Cells.Item(6, 5).Value = "  -0.06%  "
$ws.Cells.Item(7, 5).Value = "  +0.40%  "
$ws.Cells.Item(8, 5).Value = "  -0.03%  "
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0613"
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -0.16%  "
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "18.17"
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +0.25%  "
$ws.Cells.Item(11, 5).Value = "  +2.40%  "
$ws.Cells.Item(12, 4).Value = "1.816.14"
$ws.Cells.Item(12, 5).Value = "  +0.59%  "
$ws.Cells.Item(13, 4).Value = "1.578.81"
$ws.Cells.Item(13, 5).Value = "  -0.79%  "
$ws.Cells.Item(14, 5).Value = "  -0.64%  "
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.514"
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +1.24%  "
$ws.Cells.Item(16, 4).Value = "26.002.93"
$ws.Cells.Item(16, 5).Value = "  +0.45%  "
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "60.67"
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +1.18%  "
$ws.Cells.Item(18, 4).Value = "0.0₃0727"
$ws.Cells.Item(18, 5).Value = "  +0.15%  "
$ws.Cells.Item(19, 5).Value = "  -0.04%  "
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "203.91"
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +5.91%  "
$ws.Cells.Item(21, 5).Value = "  +1.20%  "
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "9.25"
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -1.25%  "
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.03"
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +1.45%  "
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.94"
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +13.68%  "
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "144.10"
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.62%  "
$ws.Cells.Item(26, 5).Value = "  -0.03%  "
$ws.Cells.Item(27, 5).Value = "  -7.46%  "
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "15.16"
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +0.52%  "
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.51"
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.91%  "
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.17"
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +0.44%  "
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0476"
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +0.75%  "
$ws.Cells.Item(32, 5).Value = "  +0.04%  "
$ws.Cells.Item(33, 5).Value = "  -4.00%  "
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.48"
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -1.09%  "
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.34"
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -0.48%  "
$ws.Cells.Item(36, 4).Value = "1.128.32"
$ws.Cells.Item(36, 5).Value = "  +2.78%  "
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0164"
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +7.58%  "
$ws.Cells.Item(38, 2).Value = "ARBITRUM"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.799"
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +2.41%  "
$ws.Cells.Item(39, 2).Value = "PaxDollar"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +0.02%  "
$ws.Cells.Item(40, 5).Value = "  -0.94%  "
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.491"
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -2.16%  "
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.780"
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -3.10%  "
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.14"
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +0.17%  "
$ws.Cells.Item(44, 4).Value = "1.726.97"
$ws.Cells.Item(44, 5).Value = "  +0.46%  "
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "92.21"
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -1.51%  "
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "53.88"
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +1.26%  "
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.49"
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -1.38%  "
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0506"
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -0.77%  "
$ws.Cells.Item(49, 5).Value = "  +0.22%  "
$ws.Cells.Item(50, 5).Value = "  +0.36%  "
$ws.Cells.Item(51, 5).Value = "  -15.70%  "
